$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- workbook-level: remove workbook protection ---
$wb.Unprotect()

# --- remove rows 5-12 (old filler data no longer present) ---
$ws.Range("A5:G12").EntireRow.Delete()

# --- column widths (best effort; engine snaps to 1/6-character grid) ---
$ws.Columns.Item(1).ColumnWidth = 9.666666666666666
$ws.Columns.Item(2).ColumnWidth = 15.833333333333334
$ws.Columns.Item(3).ColumnWidth = 14.5
$ws.Columns.Item(4).ColumnWidth = 18.333333333333332
$ws.Columns.Item(5).ColumnWidth = 15.666666666666666
$ws.Columns.Item(6).ColumnWidth = 12.0

# --- clear cells that must become empty/removed ---
$ws.Range("A3:C3").ClearContents()
$ws.Range("D3").ClearContents()
$ws.Range("F3").ClearContents()
$ws.Range("E4").ClearContents()

# --- row 1: headers ---
$ws.Range("A1").Value = "Ім'я"
$ws.Range("B1").Value = "Прізвище"
$ws.Range("C1").Value = "По батькові"
$ws.Range("D1").Value = "День народження"
$ws.Range("E1").Value = "День смерті"
$ws.Range("F1").Value = "Стать"
$ws.Range("G1").Value = "Вік"

# --- row 2 ---
$ws.Range("A2").Value = "Андрій"
$ws.Range("B2").Value = "Маслов"
$ws.Range("C2").Value = "Вікторович"
$ws.Range("D2").Value = 37002
$ws.Range("D2").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("E2").ClearContents()
$ws.Range("E2").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("F2").Value = "чоловік"
$ws.Range("G2").Value = 21

# --- row 3 ---
$ws.Range("D3").NumberFormat = "yyyy\-mm\-dd"

# --- row 4 ---
$ws.Range("A4").Value = "Юлія"
$ws.Range("B4").Value = "Павленко"
$ws.Range("C4").Value = "Семенівна"
$ws.Range("D4").Value = 36019
$ws.Range("D4").NumberFormat = "yyyy-mm-dd"
$ws.Range("F4").Value = "жінка"
$ws.Range("G4").Value = 24

# --- selection state shown in the target file ---
$ws.Range("D10").Select()
